# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value  = 1488
$ws1.Range("F11").Value = 718
$ws1.Range("F16").Value = 5774
$ws1.Range("F18").Value = 5465
$ws1.Range("F19").Value = 2062
$ws1.Range("F20").Value = 2962
$ws1.Range("F23").Value = 1663
$ws1.Range("F30").Value = 1061
$ws1.Range("F31").Value = 2239
$ws1.Range("F34").Value = 318
$ws1.Range("F35").Value = 841

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 1488
$ws4.Range("F15").Value = 718
$ws4.Range("F24").Value = 5774
$ws4.Range("F26").Value = 5465
$ws4.Range("F27").Value = 2063
$ws4.Range("F28").Value = 2962
$ws4.Range("F34").Value = 1663
$ws4.Range("F42").Value = 2239
$ws4.Range("F45").Value = 318
$ws4.Range("F46").Value = 841
